$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.989.03'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '3.753.92'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.21%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '601.98'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.12%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '165.30'
$c.ClearFormats()
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("D7").Value = '3.754.99'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.04%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.538'
$c.ClearFormats()
$ws.Range("E9").Value = '  +0.45%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.173'
$c.ClearFormats()
$ws.Range("E10").Value = '  +4.97%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.41'
$c.ClearFormats()
$ws.Range("E11").Value = '  +0.84%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.458'
$c.ClearFormats()
$ws.Range("E12").Value = '  -0.98%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '37.67'
$c.ClearFormats()
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = '4.388.28'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '3.759.30'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '69.069.04'
$ws.Range("E17").Value = '  +0.49%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '7.42'
$c.ClearFormats()
$ws.Range("E18").Value = '  +1.52%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.64'
$c.ClearFormats()
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("E20").Value = '  -0.86%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.33'
$c.ClearFormats()
$ws.Range("E21").Value = '  +5.25%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '491.02'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.93%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.724'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '84.83'
$c.ClearFormats()
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.0000148'
$c.ClearFormats()
$ws.Range("E25").Value = '  +1.17%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.ClearFormats()
$ws.Range("E26").Value = '  -1.79%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.27'
$c.ClearFormats()
$ws.Range("E27").Value = '  -1.25%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.05'
$c.ClearFormats()
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.69%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.18'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.88%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.ClearFormats()
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("B33").Value = 'WrappedeETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D33").Value = '3.906.12'
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '31.65'
$c.ClearFormats()
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").Value = '3.697.23'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -0.68%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.92'
$c.ClearFormats()
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("E38").Value = '  -0.58%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.138'
$c.ClearFormats()
$ws.Range("E39").Value = '  +4.30%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.04%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.327'
$c.ClearFormats()
$ws.Range("E41").Value = '  +0.14%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.10'
$c.ClearFormats()
$ws.Range("E42").Value = '  +8.38%  '
$ws.Range("E43").Value = '  +0.91%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '48.59'
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.88%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '424.72'
$c.ClearFormats()
$ws.Range("E45").Value = '  -2.94%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '8.43'
$c.ClearFormats()
$ws.Range("E46").Value = '  -0.96%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '40.33'
$c.ClearFormats()
$ws.Range("E48").Value = '  -0.78%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '141.52'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.788.28'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.29'
$c.ClearFormats()
$ws.Range("E51").Value = '  +7.01%  '
